$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$win = $excel.ActiveWindow
Write-Host "Window type:" $win.GetType().FullName
Write-Host "SplitRow:" $win.SplitRow
Write-Host "FreezePanes:" $win.FreezePanes
$win.ScrollRow = 5
Write-Host "ScrollRow after set:" $win.ScrollRow
